$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column Q (2020 data) mirroring the styling of column P for each row.
# Copying the P column's formatted cells first means Q inherits the exact same
# effective formatting (number format, font, borders, alignment) as column P.
$ws.Range("P4:P14").Copy($ws.Range("Q4:Q14"))

$ws.Range("Q4").Value = 2020
$ws.Range("Q5").Value = 4.5999999999999996
$ws.Range("Q6").Value = 4.2
$ws.Range("Q7").Value = 1.3
$ws.Range("Q8").Value = 10.8
$ws.Range("Q9").Value = 6.5
$ws.Range("Q10").Value = 2.9
$ws.Range("Q11").Value = 2.6
$ws.Range("Q12").Value = 13.1
$ws.Range("Q13").Value = 1
$ws.Range("Q14").Value = 1.3

# Move the active selection to column T (whole column), as recorded in the sheet view.
$ws.Range("T1:T1048576").Select()
